# GreenLeaf Schedule - "Revisione Consegna Finale WBS"
# Applies:
#   - Updates the "Fine" date of the overall "Green Leaf" summary row (F2)
#   - Adds a new WBS row 1.9 "Consegna Finale" (row 44) with its data
#   - Mirrors the view's zoom/selection change

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update end date of the top-level "Green Leaf" row (F2): 2023-01-26 -> 2023-01-13
$ws.Range("F2").Value = 44939

# 2. Copy the formatting of the previous "category" row (row 41, a bold
# "WBS summary" style row) onto the new row 44 so it matches the look of
# the other top-level WBS rows (bold font, borders, centered date format).
$ws.Range("A41:G41").Copy()
$ws.Range("A44:G44").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# 3. Populate the new row 44 (WBS 1.9 - Consegna Finale).
# The WBS label "1.9" looks like a number, so it is written through a
# helper formula cell and pasted back as a value; this keeps it stored as
# text (matching the rest of the WBS column) without altering the cell's
# number format / style that was just applied above.
$ws.Range("ZZ1").Formula = '="1.9"'
$ws.Range("ZZ1").Copy()
$ws.Range("A44").PasteSpecial(-4163)  # xlPasteValues
$ws.Application.CutCopyMode = $false
$ws.Range("ZZ1").Clear()

$ws.Range("B44").Value = "Consegna Finale"
$ws.Range("C44").Value = "1 g"
$ws.Range("D44").Value = "0h"
$ws.Range("E44").Value = 44970
$ws.Range("F44").Value = 44970
$ws.Range("G44").Value = "2;14;21;25"

# 4. Refresh the view (zoom level / selected cell) as left by the editor.
$ws.Application.ActiveWindow.Zoom = 80
$ws.Range("L31").Select()
